$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3040.1304
$ws.Range("I98").Value = 3285.85
$ws.Range("K98").Value = 3285.85
$ws.Range("M98").Value = -1787.85
$ws.Range("H116").Value = 3519.7273
$ws.Range("I116").Value = 2673
$ws.Range("J116").Value = 5001.5
$ws.Range("K116").Value = 2673
$ws.Range("L116").Value = 5001.5
$ws.Range("M116").Value = 769
$ws.Range("N116").Value = -11885.5
$ws.Range("H122").Value = 3040.1304
$ws.Range("I122").Value = 3285.85
$ws.Range("K122").Value = 9857.549999999999
$ws.Range("M122").Value = -7407.549999999999
$ws.Range("H135").Value = 41667380
$ws.Range("I135").Value = 643.4545000000001
$ws.Range("K135").Value = 5791.0905
$ws.Range("M135").Value = -3256.0905
$ws.Range("H137").Value = 1687.95
$ws.Range("I137").Value = 1558.0416
$ws.Range("J137").Value = 1882.8125
$ws.Range("K137").Value = 4674.1248
$ws.Range("L137").Value = 5648.4375
$ws.Range("M137").Value = -2124.1248
$ws.Range("N137").Value = -10748.4375
$ws.Range("H140").Value = 41239.332
$ws.Range("J140").Value = 41239.332
$ws.Range("L140").Value = 41239.332
$ws.Range("N140").Value = -51599.332
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6169.424
$ws.Range("I32").Value = 5034.36
$ws.Range("J32").Value = 11177.059
$ws.Range("K32").Value = 5034.36
$ws.Range("L32").Value = 11177.059
$ws.Range("M32").Value = -4747.36
$ws.Range("N32").Value = -11751.059
$ws.Range("H38").Value = 8061.2
$ws.Range("I38").Value = 5121.5
$ws.Range("K38").Value = 5121.5
$ws.Range("M38").Value = -4654.5
$ws.Range("H45").Value = 1115.0555
$ws.Range("I45").Value = 1047.9286
$ws.Range("K45").Value = 1047.9286
$ws.Range("M45").Value = -670.9286
$ws.Range("H61").Value = 47620316
$ws.Range("I61").Value = 58824570
$ws.Range("J61").Value = 2228.5
$ws.Range("K61").Value = 58824570
$ws.Range("L61").Value = 2228.5
$ws.Range("M61").Value = -58824358
$ws.Range("N61").Value = -2652.5
$ws.Range("H74").Value = 2147.682
$ws.Range("I74").Value = 1343.4166
$ws.Range("K74").Value = 1343.4166
$ws.Range("M74").Value = -469.4166
$ws.Range("H77").Value = 2147.682
$ws.Range("I77").Value = 1343.4166
$ws.Range("K77").Value = 6717.083000000001
$ws.Range("M77").Value = -2349.083000000001
$ws.Range("H132").Value = 3077.8064
$ws.Range("I132").Value = 3104.8572
$ws.Range("J132").Value = 3055.5293
$ws.Range("K132").Value = 9314.571599999999
$ws.Range("L132").Value = 9166.5879
$ws.Range("M132").Value = -6784.571599999999
$ws.Range("N132").Value = -14226.5879
$ws.Range("H136").Value = 47620316
$ws.Range("I136").Value = 58824570
$ws.Range("J136").Value = 2228.5
$ws.Range("K136").Value = 176473710
$ws.Range("L136").Value = 6685.5
$ws.Range("M136").Value = -176471160
$ws.Range("N136").Value = -11785.5
$ws.Range("H138").Value = 48181
$ws.Range("J138").Value = 48181
$ws.Range("L138").Value = 48181
$ws.Range("N138").Value = -58461
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7812880
$ws.Range("I94").Value = 8928870
$ws.Range("J94").Value = 952
$ws.Range("K94").Value = 8928870
$ws.Range("L94").Value = 952
$ws.Range("M94").Value = -8928419
$ws.Range("N94").Value = -1854
$ws.Range("H134").Value = 1361.3636
$ws.Range("I134").Value = 1197.5
$ws.Range("K134").Value = 3592.5
$ws.Range("M134").Value = -1057.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 52632670
$ws.Range("J16").Value = 1125.5714
$ws.Range("L16").Value = 1125.5714
$ws.Range("N16").Value = -1699.5714
$ws.Range("H22").Value = 39215.223
$ws.Range("I22").Value = 185.85715
$ws.Range("K22").Value = 185.85715
$ws.Range("M22").Value = 164.14285
$ws.Range("H31").Value = 1542.3959
$ws.Range("I31").Value = 1383.3695
$ws.Range("K31").Value = 1383.3695
$ws.Range("M31").Value = -1088.3695
$ws.Range("H34").Value = 1542.3959
$ws.Range("I34").Value = 1383.3695
$ws.Range("K34").Value = 1383.3695
$ws.Range("M34").Value = -1181.3695
$ws.Range("H92").Value = 57550.5
$ws.Range("J92").Value = 57550.5
$ws.Range("L92").Value = 57550.5
$ws.Range("N92").Value = -62542.5
$ws.Range("H113").Value = 52632670
$ws.Range("J113").Value = 1125.5714
$ws.Range("L113").Value = 1125.5714
$ws.Range("N113").Value = -5465.5714
$ws.Range("H132").Value = 3062.6428
$ws.Range("I132").Value = 2653.2222
$ws.Range("J132").Value = 3799.6
$ws.Range("K132").Value = 7959.6666
$ws.Range("L132").Value = 11398.8
$ws.Range("M132").Value = -5429.6666
$ws.Range("N132").Value = -16458.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 296
$ws.Range("I11").Value = 320
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 960
$ws.Range("L11").Value = 600
$ws.Range("M11").Value = -820
$ws.Range("N11").Value = -880
$ws.Range("H38").Value = 100.57143
$ws.Range("J38").Value = 240
$ws.Range("L38").Value = 720
$ws.Range("N38").Value = -1414
$ws.Range("H110").Value = 9574.5625
$ws.Range("I110").Value = 2002.5
$ws.Range("J110").Value = 10656.286
$ws.Range("K110").Value = 6007.5
$ws.Range("L110").Value = 31968.858
$ws.Range("M110").Value = -1917.5
$ws.Range("N110").Value = -40148.858
$ws.Range("H122").Value = 1927.1538
$ws.Range("I122").Value = 909.8
$ws.Range("J122").Value = 2563
$ws.Range("K122").Value = 8188.2
$ws.Range("L122").Value = 23067
$ws.Range("M122").Value = -5738.2
$ws.Range("N122").Value = -27967
$ws.Range("H131").Value = 21278176
$ws.Range("I131").Value = 83333790
$ws.Range("J131").Value = 1965.3715
$ws.Range("K131").Value = 250001370
$ws.Range("L131").Value = 5896.1145
$ws.Range("M131").Value = -249996330
$ws.Range("N131").Value = -15976.1145
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1846.25
$ws.Range("J126").Value = 1970.3636
$ws.Range("L126").Value = 5911.0908
$ws.Range("N126").Value = -10851.0908
$ws.Range("H141").Value = 39747
$ws.Range("J141").Value = 39747
$ws.Range("L141").Value = 39747
$ws.Range("N141").Value = -50107
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5333
$ws.Range("I40").Value = 2559.4
$ws.Range("J40").Value = 6873.8887
$ws.Range("K40").Value = 2559.4
$ws.Range("L40").Value = 6873.8887
$ws.Range("M40").Value = -2423.4
$ws.Range("N40").Value = -7145.8887
$ws.Range("H46").Value = 5610
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 6177.778
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 6177.778
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -6553.778
$ws.Range("H122").Value = 20835062
$ws.Range("I122").Value = 25001484
$ws.Range("K122").Value = 75004452
$ws.Range("M122").Value = -75002002
$ws.Range("H136").Value = 2091.6365
$ws.Range("I136").Value = 1572.5714
$ws.Range("K136").Value = 4717.7142
$ws.Range("M136").Value = -2167.7142
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10871413
$ws.Range("I122").Value = 11906666
$ws.Range("J122").Value = 1252.5
$ws.Range("K122").Value = 35719998
$ws.Range("L122").Value = 3757.5
$ws.Range("M122").Value = -35717548
$ws.Range("N122").Value = -8657.5
$ws.Range("H126").Value = 40000884
$ws.Range("I126").Value = 58823988
$ws.Range("J126").Value = 1789.25
$ws.Range("K126").Value = 176471964
$ws.Range("L126").Value = 5367.75
$ws.Range("M126").Value = -176469494
$ws.Range("N126").Value = -10307.75
$ws.Range("H137").Value = 32905
$ws.Range("J137").Value = 32905
$ws.Range("L137").Value = 32905
$ws.Range("N137").Value = -43105
